$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '42.244.36'
$ws.Cells.Item(2, 5).NumberFormat = '@'
$ws.Cells.Item(2, 5).Value = '  -0.78%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.246.46'
$ws.Cells.Item(3, 5).NumberFormat = '@'
$ws.Cells.Item(3, 5).Value = '  -1.51%  '

$ws.Cells.Item(4, 5).NumberFormat = '@'
$ws.Cells.Item(4, 5).Value = '  +0.17%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '247.50'
$ws.Cells.Item(5, 5).NumberFormat = '@'
$ws.Cells.Item(5, 5).Value = '  -1.53%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.631'
$ws.Cells.Item(6, 5).NumberFormat = '@'
$ws.Cells.Item(6, 5).Value = '  -0.36%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '75.75'
$ws.Cells.Item(7, 5).NumberFormat = '@'
$ws.Cells.Item(7, 5).Value = '  +5.21%  '

$ws.Cells.Item(8, 5).NumberFormat = '@'
$ws.Cells.Item(8, 5).Value = '  +0.08%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.629'
$ws.Cells.Item(9, 5).NumberFormat = '@'
$ws.Cells.Item(9, 5).Value = '  -2.15%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '40.23'
$ws.Cells.Item(10, 5).NumberFormat = '@'
$ws.Cells.Item(10, 5).Value = '  +4.00%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0951'
$ws.Cells.Item(11, 5).NumberFormat = '@'
$ws.Cells.Item(11, 5).Value = '  -1.74%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '7.24'
$ws.Cells.Item(12, 5).NumberFormat = '@'
$ws.Cells.Item(12, 5).Value = '  -1.89%  '

$ws.Cells.Item(13, 5).NumberFormat = '@'
$ws.Cells.Item(13, 5).Value = '  -1.40%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '2.587.95'
$ws.Cells.Item(14, 5).NumberFormat = '@'
$ws.Cells.Item(14, 5).Value = '  -1.19%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '14.92'
$ws.Cells.Item(15, 5).NumberFormat = '@'
$ws.Cells.Item(15, 5).Value = '  -1.05%  '

$ws.Cells.Item(16, 5).NumberFormat = '@'
$ws.Cells.Item(16, 5).Value = '  -2.04%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.238.72'
$ws.Cells.Item(17, 5).NumberFormat = '@'
$ws.Cells.Item(17, 5).Value = '  -1.34%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '42.238.33'
$ws.Cells.Item(18, 5).NumberFormat = '@'
$ws.Cells.Item(18, 5).Value = '  -0.73%  '

$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '0.0₃0980'
$ws.Cells.Item(19, 5).NumberFormat = '@'
$ws.Cells.Item(19, 5).Value = '  -1.68%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.17'
$ws.Cells.Item(20, 5).NumberFormat = '@'
$ws.Cells.Item(20, 5).Value = '  -2.35%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '71.58'
$ws.Cells.Item(21, 5).NumberFormat = '@'
$ws.Cells.Item(21, 5).Value = '  -1.41%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '231.90'
$ws.Cells.Item(22, 5).NumberFormat = '@'
$ws.Cells.Item(22, 5).Value = '  -1.46%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '2.21'
$ws.Cells.Item(23, 5).NumberFormat = '@'
$ws.Cells.Item(23, 5).Value = '  -4.46%  '

$ws.Cells.Item(24, 5).NumberFormat = '@'
$ws.Cells.Item(24, 5).Value = '  -0.01%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '3.74'
$ws.Cells.Item(25, 5).NumberFormat = '@'
$ws.Cells.Item(25, 5).Value = '  -4.81%  '

$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '11.12'
$ws.Cells.Item(26, 5).NumberFormat = '@'
$ws.Cells.Item(26, 5).Value = '  -4.46%  '

$ws.Cells.Item(27, 5).NumberFormat = '@'
$ws.Cells.Item(27, 5).Value = '  -4.33%  '

$ws.Cells.Item(28, 2).NumberFormat = '@'
$ws.Cells.Item(28, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28, 3).NumberFormat = '@'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '7.08'
$ws.Cells.Item(28, 5).NumberFormat = '@'
$ws.Cells.Item(28, 5).Value = '  +9.07%  '

$ws.Cells.Item(29, 2).NumberFormat = '@'
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).NumberFormat = '@'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.16'
$ws.Cells.Item(29, 5).NumberFormat = '@'
$ws.Cells.Item(29, 5).Value = '  -1.10%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '168.67'
$ws.Cells.Item(30, 5).NumberFormat = '@'
$ws.Cells.Item(30, 5).Value = '  +0.85%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '20.57'
$ws.Cells.Item(31, 5).NumberFormat = '@'
$ws.Cells.Item(31, 5).Value = '  -2.43%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.0850'
$ws.Cells.Item(32, 5).NumberFormat = '@'
$ws.Cells.Item(32, 5).Value = '  +5.72%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '32.80'
$ws.Cells.Item(33, 5).NumberFormat = '@'
$ws.Cells.Item(33, 5).Value = '  +3.66%  '

$ws.Cells.Item(34, 5).NumberFormat = '@'
$ws.Cells.Item(34, 5).Value = '  -6.20%  '

$ws.Cells.Item(35, 5).NumberFormat = '@'
$ws.Cells.Item(35, 5).Value = '  -0.09%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '4.51'
$ws.Cells.Item(36, 5).NumberFormat = '@'
$ws.Cells.Item(36, 5).Value = '  -5.38%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.78'
$ws.Cells.Item(37, 5).NumberFormat = '@'
$ws.Cells.Item(37, 5).Value = '  +0.70%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0298'
$ws.Cells.Item(38, 5).NumberFormat = '@'
$ws.Cells.Item(38, 5).Value = '  -3.74%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '13.16'
$ws.Cells.Item(39, 5).NumberFormat = '@'
$ws.Cells.Item(39, 5).Value = '  -6.26%  '

$ws.Cells.Item(40, 2).NumberFormat = '@'
$ws.Cells.Item(40, 2).Value = 'THORChain'
$ws.Cells.Item(40, 3).NumberFormat = '@'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '5.98'
$ws.Cells.Item(40, 5).NumberFormat = '@'
$ws.Cells.Item(40, 5).Value = '  -0.22%  '

$ws.Cells.Item(41, 2).NumberFormat = '@'
$ws.Cells.Item(41, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(41, 3).NumberFormat = '@'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '2.21'
$ws.Cells.Item(41, 5).NumberFormat = '@'
$ws.Cells.Item(41, 5).Value = '  -5.68%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '117.74'
$ws.Cells.Item(42, 5).NumberFormat = '@'
$ws.Cells.Item(42, 5).Value = '  +21.26%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.203'
$ws.Cells.Item(43, 5).NumberFormat = '@'
$ws.Cells.Item(43, 5).Value = '  -4.46%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '60.23'
$ws.Cells.Item(44, 5).NumberFormat = '@'
$ws.Cells.Item(44, 5).Value = '  -2.85%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '8.71'
$ws.Cells.Item(45, 5).NumberFormat = '@'
$ws.Cells.Item(45, 5).Value = '  -5.59%  '

$ws.Cells.Item(46, 5).NumberFormat = '@'
$ws.Cells.Item(46, 5).Value = '  -2.62%  '

$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.998'
$ws.Cells.Item(47, 5).NumberFormat = '@'
$ws.Cells.Item(47, 5).Value = '  -0.30%  '

$ws.Cells.Item(48, 5).NumberFormat = '@'
$ws.Cells.Item(48, 5).Value = '  -4.09%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.17'
$ws.Cells.Item(49, 5).NumberFormat = '@'
$ws.Cells.Item(49, 5).Value = '  -1.90%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '4.26'
$ws.Cells.Item(50, 5).NumberFormat = '@'
$ws.Cells.Item(50, 5).Value = '  -12.90%  '

$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '4.14'
$ws.Cells.Item(51, 5).NumberFormat = '@'
$ws.Cells.Item(51, 5).Value = '  -2.54%  '
